$wb = $excel.ActiveWorkbook

# --- 1. "总计" sheet: insert a new row 2 holding the 2022-Q4 summary ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
# Copy formatting from the row below (now row 3, still carrying the same
# style as the rest of the index column) onto the freshly inserted row.
$ws1.Range("A3:D3").Copy()
$ws1.Range("A2:D2").PasteSpecial(-4122)
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 2

# --- 2. Insert the brand-new "2022-Q4" worksheet right after "总计" ---
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$newSheet.Name = "2022-Q4"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows (fund code / numeric-looking figures are stored as text in the
# source data, so they're entered with a leading apostrophe to stop Excel
# from auto-coercing them into numbers, e.g. losing the leading zero of
# "010363").
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'010363"
$newSheet.Range("C2").Value = "信澳匠心臻选两年持有期混合"
$newSheet.Range("D2").Value = "'37.44"
$newSheet.Range("E2").Value = "'92.71"
$newSheet.Range("F2").Value = "'2.81"
$newSheet.Range("G2").Value = "'1.0521"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'010963"
$newSheet.Range("C3").Value = "信澳周期动力混合A"
$newSheet.Range("D3").Value = "'24.33"
$newSheet.Range("E3").Value = "'92.85"
$newSheet.Range("F3").Value = "'3.21"
$newSheet.Range("G3").Value = "'0.7810"
$newSheet.Range("H3").Value = 8

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'015455"
$newSheet.Range("C4").Value = "信澳周期动力混合C"
$newSheet.Range("D4").Value = "'5.15"
$newSheet.Range("E4").Value = "'92.85"
$newSheet.Range("F4").Value = "'3.21"
$newSheet.Range("G4").Value = "'0.1653"
$newSheet.Range("H4").Value = 8

# Writing text that looks numeric auto-applies a "Text" number format/style;
# strip that back off so these cells fall back to the default (unstyled)
# look, matching the other data cells in the workbook.
$newSheet.Range("B2:G4").Style = "Normal"

# Re-apply the real visual style used throughout this workbook (bold+border
# header row, bold+border index column) by copying it from the sibling
# "2022-Q3" sheet, which already carries the exact same formatting.
$srcSheet = $wb.Worksheets.Item(3)
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2:A4").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# Restore the original active sheet/selection ("总计", cell A1) since adding
# a sheet otherwise leaves the new sheet focused.
$null = $newSheet.Range("A1").Select()
$ws1.Activate()
$null = $ws1.Range("A1").Select()
